# Update generation-mix figures (MORE Trading Node file path change).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")

$ws.Range("A2").Value = 30928.199
$ws.Range("C2").Value = 16500

$ws.Range("A3").Value = 29619.8635
$ws.Range("B3").Value = 10000
$ws.Range("C3").Value = 13500
$ws.Range("D3").Value = 20000
$ws.Range("E3").Value = 10000
$ws.Range("F3").Value = 10000
$ws.Range("G3").Value = 20000
$ws.Range("J3").Value = 10000

$ws.Range("A4").Value = 27793.1675
$ws.Range("G4").Value = 20000

$ws.Range("A5").Value = 26218.02
$ws.Range("G5").Value = 20000

$ws.Range("A6").Value = 25305.1185

$ws.Range("A7").Value = 25440.605
$ws.Range("G7").Value = 0

$ws.Range("A8").Value = 25564.52277227723
$ws.Range("G8").Value = 0

$ws.Range("A9").Value = 26075.5535
$ws.Range("G9").Value = 0

$ws.Range("A10").Value = 28306.4105
$ws.Range("G10").Value = 0

$ws.Range("A11").Value = 26016.9995
$ws.Range("G11").Value = 10000

$ws.Range("A12").Value = 24997.2015
$ws.Range("G12").Value = 10000

$ws.Range("A13").Value = 24407.995
$ws.Range("G13").Value = 10000

$ws.Range("A14").Value = 23983.6615
$ws.Range("G14").Value = 10000

$ws.Range("A15").Value = 25520.364
$ws.Range("C15").Value = 12500
$ws.Range("D15").Value = 10000
$ws.Range("E15").Value = 5000
$ws.Range("F15").Value = 5000

$ws.Range("A16").Value = 26075.793

$ws.Range("A17").Value = 25427.0925

$ws.Range("A18").Value = 23201.19782374277

$ws.Range("A19").Value = 20688.49828448191

$ws.Range("C24").Value = 22000

$ws.Range("C25").Value = 12500
